$p = $ppt.ActivePresentation

# -----------------------------------------------------------------
# 1) Slide 2: remove the "Cloud 50" shape and its "Elbow Connector 51"
#    connector (the connector references the cloud, so delete it first).
# -----------------------------------------------------------------
$s2 = $p.Slides.Item(2)

foreach ($shapeName in @("Elbow Connector 51", "Cloud 50")) {
    for ($i = $s2.Shapes.Count; $i -ge 1; $i--) {
        $shp = $s2.Shapes.Item($i)
        if ($shp.Name -eq $shapeName) {
            $shp.Delete()
        }
    }
}

# -----------------------------------------------------------------
# 2) Refresh the "updates automatically" date footer field from
#    10/28/2016 to 11/3/2016 everywhere it is cached: the slide
#    master, every slide layout, and the notes master.
# -----------------------------------------------------------------
function Update-DateShapes($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            if ($shp.TextFrame.TextRange.Text -eq "10/28/2016") {
                $shp.TextFrame.TextRange.Text = "11/3/2016"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateShapes $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DateShapes $layout.Shapes
}

$notesMaster = $p.NotesMaster
Update-DateShapes $notesMaster.Shapes
